$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.163.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.779.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.54%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.552"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0657"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.033.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.767.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "

$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.149.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0739"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.53%  "

$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("E28").Value = "  -1.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  -4.73%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.438.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.20%  "

$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.624"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "

$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.889"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.32%  "

$ws.Range("E43").Value = "  -5.56%  "

$ws.Range("E44").Value = "  -2.73%  "

$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.935.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.30%  "
